# Auto-generated: apply CTA violent-crime YTD data refresh (2025-11-23)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('F2').Value = 87
$ws.Range('J2').Value = 115
$ws.Range('L2').Value = 123
$ws.Range('D3').Value = 130
$ws.Range('E3').Value = 139
$ws.Range('F3').Value = 131
$ws.Range('I3').Value = 186
$ws.Range('K3').Value = 212
$ws.Range('C6').Value = 460
$ws.Range('D6').Value = 394
$ws.Range('F6').Value = 499
$ws.Range('G6').Value = 425
$ws.Range('I6').Value = 485
$ws.Range('K6').Value = 486
$ws.Range('L6').Value = 418
$ws.Range('C7').Value = 611
$ws.Range('D7').Value = 620
$ws.Range('E7').Value = 666
$ws.Range('F7').Value = 725
$ws.Range('G7').Value = 651
$ws.Range('I7').Value = 811
$ws.Range('J7').Value = 755
$ws.Range('K7').Value = 861
$ws.Range('L7').Value = 803

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L6').Value = 32
$ws.Range('L7').Value = 44

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('D5').Value = 10
$ws.Range('D6').Value = 13

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('G5').Value = 5
$ws.Range('G6').Value = 10

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('D4').Value = 4
$ws.Range('D5').Value = 5

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('F2').Value = 2
$ws.Range('K4').Value = 4
$ws.Range('F5').Value = 12
$ws.Range('K5').Value = 9

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I19').Value = 22
$ws.Range('K29').Value = 22
$ws.Range('L32').Value = 44
$ws.Range('G41').Value = 10
$ws.Range('E47').Value = 17
$ws.Range('D50').Value = 13
$ws.Range('I52').Value = 7
$ws.Range('C53').Value = 53
$ws.Range('F68').Value = 4
$ws.Range('D75').Value = 5
$ws.Range('K76').Value = 28
$ws.Range('F77').Value = 19
$ws.Range('I77').Value = 47
$ws.Range('K79').Value = 8
$ws.Range('L79').Value = 11
$ws.Range('D80').Value = 5
$ws.Range('F82').Value = 12
$ws.Range('K82').Value = 9
$ws.Range('I85').Value = 5
$ws.Range('J92').Value = 6
$ws.Range('C98').Value = 611
$ws.Range('D98').Value = 620
$ws.Range('E98').Value = 666
$ws.Range('F98').Value = 725
$ws.Range('G98').Value = 651
$ws.Range('I98').Value = 811
$ws.Range('J98').Value = 755
$ws.Range('K98').Value = 861
$ws.Range('L98').Value = 803

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('C6').Value = 35
$ws.Range('C7').Value = 53

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J2').Value = 1
$ws.Range('J5').Value = 6

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range('L2').Value = 2
$ws.Range('K3').Value = 2
$ws.Range('K6').Value = 8
$ws.Range('L6').Value = 11

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K5').Value = 13
$ws.Range('K6').Value = 28

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('K5').Value = 17
$ws.Range('K6').Value = 22

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('F3').Value = 7
$ws.Range('I3').Value = 9
$ws.Range('F7').Value = 19
$ws.Range('I7').Value = 47

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I6').Value = 14
$ws.Range('I7').Value = 22

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('I4').Value = 2
$ws.Range('I5').Value = 5

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('I5').Value = 5
$ws.Range('I6').Value = 7

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('E3').Value = 3
$ws.Range('E6').Value = 17

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('D3').Value = 3
$ws.Range('D5').Value = 5

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range('F6').Value = 3
$ws.Range('F7').Value = 4
